$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 6
$ws1.Range("F15").Value = 12921
$ws1.Range("G15").Value = 60
$ws1.Range("F16").Value = 10
$ws1.Range("F17").Value = 5303

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 144

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 6
$ws4.Range("F15").Value = 12921
$ws4.Range("G15").Value = 60
$ws4.Range("F16").Value = 144
$ws4.Range("F18").Value = 10
$ws4.Range("F19").Value = 5303
